$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The handback transform failed for the 85874fdf... file in both target
# languages, so its "Status" flips from "Ready for handoff" to
# "Handback transform failed" everywhere it's reported: the per-language
# detail sheets and the rollup Overview sheet.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Populate the "Error Detail" column (P) for that same row on both the
# zh-cn and de-de target sheets with the handback-transform error message.
$zhcn.Range("P3").Value = "Handback file name: ixktxe0y.mo2 is different with handoff file name: 85874fdf-7ff8-429d-b7b8-9e6b73751bed.162ef1162fb13e280d15f927cddcc4ed20ff3ba2.zh-cn."
$dede.Range("P3").Value = "Handback file name: ixktxe0y.mo2 is different with handoff file name: 85874fdf-7ff8-429d-b7b8-9e6b73751bed.162ef1162fb13e280d15f927cddcc4ed20ff3ba2.de-de."

# Widen the now much-longer "Error Detail" column on both sheets to fit
# the new message text (matches the other 40-wide text columns).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
